$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the duplicate "Barack Obama" row (row 47 is an exact duplicate of row 46).
#    A normal row delete shifts rows 48/49 up by one, leaving the trailing placeholder
#    row (old row 49) as new row 48.
$ws.Rows(47).Delete()

# 2. Insert a new column before the current column D ("prior"), to host a cleaned-up
#    "president-fixed" column built from PROPER(president).
$ws.Columns("D").Insert()

# 3. Header for the new column.
$ws.Range("D1").Value = "president-fixed"

# 4. Formula column: D2 gets its own formula, D3:D47 fill down as a shared formula
#    (mirrors entering =PROPER(C2) in D2 and dragging the fill handle to D47).
$ws.Range("D2").Formula = "=PROPER(C2)"
$ws.Range("D3:D47").Formula = "=PROPER(C3)"

# 5. Re-materialize the trailing blank placeholder row that the delete above pushed out
#    (old row 49 had a lone empty formatted cell; keep that pattern two rows down in H).
$ws.Range("H49").Style = "Normal"

# 6. Update the view: scroll back to the top and select D1 (matches the new selection,
#    no more frozen/scrolled topLeftCell).
$ws.Range("D1").Select()
